$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personas")

# New row for a fresh "Titular" (primary) occupant, folio 1002.
# Some values look numeric/date-like but must be stored as literal text
# (matching the source system's inlineStr export), so we temporarily force
# a text number format before assigning them, then restore the default
# "Normal" style so the written cell ends up with no explicit style index -
# same as the rest of the sheet's data rows.

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "1002"
$ws.Range("A5").Style = "Normal"

$ws.Range("B5").Value = "juan perez"
$ws.Range("C5").Value = "8sdausd8a0sd8j"
$ws.Range("D5").Value = 27

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1998-01-14"
$ws.Range("E5").Style = "Normal"

$ws.Range("F5").Value = "americana"
$ws.Range("G5").Value = "masculino"
$ws.Range("H5").Value = "Titular"

# Empty text cells (tutor_folio, fecha_salida, motivo_salida are blank for
# this record) - a bare apostrophe yields an empty *text* cell instead of a
# truly blank/number cell, then the style reset drops the quote-prefix style.
$ws.Range("I5").Value = "'"
$ws.Range("I5").Style = "Normal"

$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = "2025-12-17 13:27:20"
$ws.Range("J5").Style = "Normal"

$ws.Range("K5").Value = 0

$ws.Range("L5").Value = "'"
$ws.Range("L5").Style = "Normal"

$ws.Range("M5").Value = "'"
$ws.Range("M5").Style = "Normal"
